$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in row 4 (P4 and Q4) ---
$ws.Range("P4").Value2 = 0.09130340807234763
$ws.Range("Q4").Value2 = 0.0748624809945284

# --- Add new column R: header year 2021 in R3 (same look as O3/P3/Q3) ---
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value2 = 2021

# --- Add new column R: data value in R4 (same look as P4/Q4, with its own font) ---
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Font.Name = "Times New Roman"
$ws.Range("R4").Font.Size = 9
$ws.Range("R4").Value2 = 0.06446742133754044

# --- Update selection / active cell on the sheet view ---
$ws.Range("O10").Select()

$wb.Save()
